$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of column M (the column to the left of the insertion
# point) so the newly-inserted column N can inherit it, mirroring Excel's
# native "Insert Column" behaviour (format/width copied from the left).
$mWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N; everything from N onward
# (old N/O/P) shifts right to O/P/Q.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, and update its selection.
$ws.Activate() | Out-Null
$ws.Range("K14").Select() | Out-Null
